$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")
# capture the existing style object reference from C16 before changing value
$style = $ws.Range("C16").Style
Write-Host "style name:" $style
$ws.Range("C16").Value = 48
$ws.Range("C16").Style = $style
